$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (prices) are stored as literal text so
# formatted numbers like "1.000" or "30.693.72" are preserved exactly.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.693.72'
$ws.Range("E2").Value = '  +1.66%  '

$ws.Range("D3").Value = '1.895.29'
$ws.Range("E3").Value = '  +2.28%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '238.69'
$ws.Range("E5").Value = '  +1.44%  '

$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").Value = '0.4837'
$ws.Range("E7").Value = '  +1.36%  '

$ws.Range("D8").Value = '0.2889'
$ws.Range("E8").Value = '  +3.04%  '

$ws.Range("D9").Value = '0.06563'
$ws.Range("E9").Value = '  +2.08%  '

$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.880.77'
$ws.Range("E10").Value = '  +1.58%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '16.90'
$ws.Range("E11").Value = '  +4.50%  '

$ws.Range("D12").Value = '0.07466'
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").Value = '5.122'
$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("D14").Value = '88.14'
$ws.Range("E14").Value = '  +1.53%  '

$ws.Range("D15").Value = '0.6702'
$ws.Range("E15").Value = '  +4.13%  '

$ws.Range("D16").Value = '30.682.86'
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("D17").Value = '13.28'
$ws.Range("E17").Value = '  +1.49%  '

$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  -0.12%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007592'
$ws.Range("E19").Value = '  +0.88%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '231.27'
$ws.Range("E20").Value = '  +3.49%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.088.60'
$ws.Range("E21").Value = '  -0.25%  '

$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.13%  '

$ws.Range("D23").Value = '5.284'
$ws.Range("E23").Value = '  +0.69%  '

$ws.Range("D24").Value = '6.219'
$ws.Range("E24").Value = '  +2.63%  '

$ws.Range("D25").Value = '170.02'
$ws.Range("E25").Value = '  +4.61%  '

$ws.Range("D26").Value = '9.385'
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("D27").Value = '18.90'
$ws.Range("E27").Value = '  +2.50%  '

$ws.Range("D28").Value = '1.966'
$ws.Range("E28").Value = '  +3.03%  '

$ws.Range("D29").Value = '0.1029'
$ws.Range("E29").Value = '  +12.44%  '

$ws.Range("E30").Value = '  -3.20%  '

$ws.Range("D31").Value = '4.348'
$ws.Range("E31").Value = '  +3.34%  '

$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  +2.59%  '

$ws.Range("D33").Value = '0.05069'
$ws.Range("E33").Value = '  +2.44%  '

$ws.Range("D34").Value = '1.216'
$ws.Range("E34").Value = '  +7.10%  '

$ws.Range("D35").Value = '0.7556'
$ws.Range("E35").Value = '  +4.56%  '

$ws.Range("D36").Value = '0.9995'
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.87%  '

$ws.Range("D38").Value = '0.01880'
$ws.Range("E38").Value = '  +2.87%  '

$ws.Range("D39").Value = '2.651'
$ws.Range("E39").Value = '  +2.18%  '

$ws.Range("D40").Value = '0.9210'
$ws.Range("E40").Value = '  +2.68%  '

$ws.Range("D41").Value = '2.073'
$ws.Range("E41").Value = '  +2.47%  '

$ws.Range("D42").Value = '107.10'
$ws.Range("E42").Value = '  +1.42%  '

$ws.Range("D43").Value = '0.4309'
$ws.Range("E43").Value = '  +2.24%  '

$ws.Range("D45").Value = '5.674'
$ws.Range("E45").Value = '  -3.66%  '

$ws.Range("D46").Value = '7.445'
$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("D47").Value = '64.30'
$ws.Range("E47").Value = '  +0.96%  '

$ws.Range("D48").Value = '0.1279'
$ws.Range("E48").Value = '  -1.97%  '

$ws.Range("D49").Value = '1.497'
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").Value = '8.971'
$ws.Range("E50").Value = '  +3.89%  '

$ws.Range("D51").Value = '34.08'
$ws.Range("E51").Value = '  +1.41%  '

